## Apply the "Add files via upload" change to TC_PartnerManagement_Master.xlsx
## - Adds a new "Approval" feature block (7 test cases, rows 105-111) to the Master sheet
## - Widens column A slightly (new col split) and updates the used-range dimension
## - Moves the sheet selection to the header row (A1:G1)
## - Reassigns the "Credential Services" drop-down validation to cover the new rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New testcase rows (105-111) - "Approval" feature under Partner Mgmt
#    Fill column-by-column (B,C,D,E,F) for every row first, then go back and
#    fill column A (TestCase_No) - this mirrors how the rows were authored
#    and keeps the shared-string table ordering stable.
# ---------------------------------------------------------------------------
$newRows = @(
  @{ B = "Partner Mgmt"; C = "Approval"; D = "Functional"; E = "Verify new partner with new policy request for api key"; F = "It should be manual approval" },
  @{ B = "Partner Mgmt"; C = "Approval"; D = "Functional"; E = "Verify same partner with same policy request for api key for subsequent times"; F = "it should be auto approval" },
  @{ B = "Partner Mgmt"; C = "Approval"; D = "Functional"; E = "Verify same partner with different policy request for api key"; F = "It should be manual approval" },
  @{ B = "Partner Mgmt"; C = "Approval"; D = "Functional"; E = "Verify same partner with same policy  request for api key but policy deactivated "; F = "It should not get auto approval" },
  @{ B = "Partner Mgmt"; C = "Approval"; D = "Functional"; E = "Verify same partner with same policy  request for api key   but partner  deactivated"; F = "It should not get auto approval" },
  @{ B = "Partner Mgmt"; C = "Approval"; D = "Functional"; E = "Verify same partner with same policy request for api key for subsequent times and try to do manual approval"; F = "it should get already approved" },
  @{ B = "Partner Mgmt"; C = "Approval"; D = "Functional"; E = "Verify auto approval if previous api key's is deactive"; F = "It should not get auto approval" }
)

$startRow = 105
for ($i = 0; $i -lt $newRows.Count; $i++) {
  $r = $startRow + $i
  $data = $newRows[$i]
  $ws.Cells.Item($r, 2).Value = $data.B
  $ws.Cells.Item($r, 3).Value = $data.C
  $ws.Cells.Item($r, 4).Value = $data.D
  $ws.Cells.Item($r, 5).Value = $data.E
  $ws.Cells.Item($r, 6).Value = $data.F
  $ws.Rows($r).RowHeight = 15
}

$testCaseIds = @("PM_Approval_01", "PM_Approval_02", "PM_Approval_03", "PM_Approval_04", "PM_Approval_05", "PM_Approval_06", "PM_Approval_07")
for ($i = 0; $i -lt $testCaseIds.Count; $i++) {
  $ws.Cells.Item($startRow + $i, 1).Value = $testCaseIds[$i]
}

# Wrap text to match the rest of the sheet's look (TestCase_No / Type / Scenario columns,
# plus the Expected Result column on the very last new row)
$ws.Range("A105:A111").WrapText = $true
$ws.Range("D105:D111").WrapText = $true
$ws.Range("E105:E111").WrapText = $true
$ws.Range("F111").WrapText = $true

# ---------------------------------------------------------------------------
# 2. Column widths - column A now gets its own (slightly wider) width instead
#    of sharing the A:B group width.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 15.166666666666666

# ---------------------------------------------------------------------------
# 3. Data validation - the "Credential Services" list now needs to cover the
#    new rows (83:111) instead of stopping at 104, and the plain Partner-Mgmt
#    list (B105:B272) needs to start after the new block (B112:B272).
# ---------------------------------------------------------------------------
$ws.Range("B83:B111").Validation.Delete()
$ws.Range("B83:B111").Validation.Add(3, 1, 1, """Admin, Credential Services, IDA, Partner Mgmt, Pre Registration, Registration Client, Registration Processor, Resident Services""")

# ---------------------------------------------------------------------------
# 4. Move the active selection up to the header row, like the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("A1:G1").Select() | Out-Null
